# Gastos.xlsx update — add the two new expense rows for 08/10/2022 and 13/10/2022
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tbl = $ws.ListObjects.Item("Tabla1")

# --- Row 35: 38000 / "construccion chiquero y clavos" / 08/10/2022 (serial 44842) ---
$row35 = $tbl.ListRows.Add()
$row35.Range.Item(1, 1).Value = 38000
$row35.Range.Item(1, 2).Value = "construccion chiquero y clavos"
$row35.Range.Item(1, 3).Value = 44842

# --- Row 36: 15000 / "yodo y guantes" / 13/10/2022 (serial 44847) ---
$row36 = $tbl.ListRows.Add()
$row36.Range.Item(1, 1).Value = 15000
$row36.Range.Item(1, 2).Value = "yodo y guantes"
$row36.Range.Item(1, 3).Value = 44847

# Match the date formatting already used by the rest of column C (copy format
# from the preceding date cell so both new cells reuse the same style, rather
# than creating a brand new custom numFmt).
$ws.Range("C34").Copy() | Out-Null
$ws.Range("C35").PasteSpecial(-4122) | Out-Null
$ws.Range("C34").Copy() | Out-Null
$ws.Range("C36").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# The running-total formula in G16 is outside the table, so it does not
# auto-expand with the table — extend it manually to cover the new rows.
$ws.Range("G16").Formula = "=SUM(A2:A36)"

# Move the selection to reflect where the user ended up after entering the
# new data (mirrors the author's last edit position).
$ws.Range("D35").Select()
